# Ambermoon-Advanced workbook update: "Nearly finalized cave upper level"
#
# Content changes:
#  - GlobalVars: var 226 changes meaning from "unused" to "You opened the
#    treasure room in dwarf mine"; new var 235 documented.
#  - MapChanges: cave teleport note updated (deactivated at start / black
#    tiles note added); new MapChanges row for the old dwarf mine (map 438).
#  - Active sheet / selections updated to reflect where the author was
#    working (GlobalVars active, Chests no longer active).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# GlobalVars sheet
# ---------------------------------------------------------------------
$wsGlobalVars = $wb.Worksheets.Item("GlobalVars")

# Var 226 no longer unused - it is now set when the dwarf mine treasure
# room is opened.
$wsGlobalVars.Range("A9").Value = "226: You opened the treasure room in dwarf mine"

# New var 235 appended as a new row.
$wsGlobalVars.Range("A18").Value = "235: You fell through a cave hole so that it create a hole below as well"

# ---------------------------------------------------------------------
# MapChanges sheet
# ---------------------------------------------------------------------
$wsMapChanges = $wb.Worksheets.Item("MapChanges")

# Row 2 (map 157) note text updated + row made a bit taller to fit it.
$wsMapChanges.Range("C2").Value = "Added teleport to renovated house (with condition)`nAdded wind gate teleport (with condition)`nAdded cave teleport (deactivated at start)`nMade back tiles below later cave black`nNPC Karl can create a wind gate there"
$wsMapChanges.Rows.Item(2).RowHeight = 75

# New row 8: map 438 "Old dwarf mine".
$wsMapChanges.Range("A8").Value = 438
$wsMapChanges.Range("B8").Value = "Old dwarf mine"
$wsMapChanges.Range("C8").Value = "Global var 226 is now set when you open the treasure room"

# ---------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------
# Update the selection remembered on sheets that are not the final active
# sheet first, so they don't themselves become "active" in the process.
$wsMapChanges.Range("B10").Select() | Out-Null

# Finally, make GlobalVars the active sheet/tab with its own selection -
# this also clears "tabSelected" from whichever sheet (Chests) had it.
$wsGlobalVars.Activate() | Out-Null
$wsGlobalVars.Range("A19").Select() | Out-Null

# ---------------------------------------------------------------------
# Workbook window (best effort - geometry may not round-trip in this
# sandboxed host, but set it anyway in case it does).
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 17640
